# BAOCAO VA SUA BAN LE KHUYEN MAI
# Rename the purchase-import template header labels and move the active
# cell selection on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column headers row 1: STT (A1) and Barcode (D1) stay the same;
# the remaining labels are renamed to the new terms.
$ws.Range("B1").Value2 = "Masieuthi"
$ws.Range("C1").Value2 = "Tenviettat"
$ws.Range("E1").Value2 = "Giabanbuoncovat"
$ws.Range("F1").Value2 = "Giabanlecovat"
$ws.Range("G1").Value2 = "Makhachhang"
$ws.Range("H1").Value2 = "Soluong"

# The active cell/selection moves from H1 to D8.
$ws.Range("D8").Select()
